$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(-0.04471486901827749, 1.538188180770143, 8.820941227871542, 2.970006940711005, 3.039554476939365, 22)
    3 = @(0.5702471346781629, 2.208141932403022, 19.53235185232473, 4.419542041017907, 4.490827157729617, 21)
    4 = @(-0.9928311212697511, 1.494578616830107, 8.573381211930261, 2.928033676706991, 2.826131366665581, 20)
    5 = @(0.04445486954785217, 0.5092960956497293, 0.4755638644454255, 0.689611386539858, 0.7070346899928552, 19)
    6 = @(0.125870100697262, 0.8822434314276495, 1.988968042097669, 1.410307782754413, 1.445403376982404, 18)
    7 = @(-0.115663227681736, 0.7181437739269187, 1.386566374226969, 1.177525530180543, 1.207895966749015, 17)
    8 = @(0.05244267050036469, 0.5335846545226226, 0.4934662843964752, 0.7024715541546683, 0.7250980603934585, 15)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Cells.Item($row, 2).Value = $values[0]
    $ws.Cells.Item($row, 3).Value = $values[1]
    $ws.Cells.Item($row, 4).Value = $values[2]
    $ws.Cells.Item($row, 5).Value = $values[3]
    $ws.Cells.Item($row, 6).Value = $values[4]
    $ws.Cells.Item($row, 7).Value = $values[5]
}
